$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.694.13"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.603.97"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'212.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'0.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "'29.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.95%  "
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.834.41"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "1.607.38"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("D15").Value = "29.730.32"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "'3.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "'64.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "'241.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'8.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.30%  "
$ws.Range("D20").Value = "0.0₃0704"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'9.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.63%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("D25").Value = "'157.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "'15.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").Value = "1.423.45"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").Value = "'1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.62%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'2.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +2.41%  "
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("D41").Value = "'55.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.79%  "
$ws.Range("E42").Value = "  +6.84%  "
$ws.Range("D43").Value = "'0.819"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "'67.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("D47").Value = "'0.995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.78%  "
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").Value = "1.743.62"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "'86.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  -1.71%  "
